$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.332.55'
$ws.Range('E2').Value = '  +1.33%  '

$ws.Range('D3').Value = '1.905.26'
$ws.Range('E3').Value = '  +1.19%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '322.30'
$ws.Range('E5').Value = '  -2.32%  '

$ws.Range('E6').Value = '  +0.03%  '

$ws.Range('D7').Value = '0.4721'
$ws.Range('E7').Value = '  +2.87%  '

$ws.Range('D8').Value = '0.4033'
$ws.Range('E8').Value = '  -1.09%  '

$ws.Range('E9').Value = '  +0.91%  '

$ws.Range('D10').Value = '0.9936'
$ws.Range('E10').Value = '  +0.20%  '

$ws.Range('D11').Value = '22.62'
$ws.Range('E11').Value = '  +4.72%  '

$ws.Range('D12').Value = '1.901.37'
$ws.Range('E12').Value = '  +1.74%  '

$ws.Range('D13').Value = '5.866'
$ws.Range('E13').Value = '  -0.61%  '

$ws.Range('D14').Value = '7.060'
$ws.Range('E14').Value = '  +0.02%  '

$ws.Range('D15').Value = '89.28'

$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.03%  '

$ws.Range('D17').Value = '0.06621'
$ws.Range('E17').Value = '  +0.74%  '

$ws.Range('D18').Value = '0.00001029'
$ws.Range('E18').Value = '  +0.34%  '

$ws.Range('D19').Value = '17.57'
$ws.Range('E19').Value = '  +1.22%  '

$ws.Range('E20').Value = '  -0.01%  '

$ws.Range('D21').Value = '29.329.83'
$ws.Range('E21').Value = '  +1.45%  '

$ws.Range('D22').Value = '5.514'
$ws.Range('E22').Value = '  +2.01%  '

$ws.Range('E23').Value = '  -0.24%  '

$ws.Range('D24').Value = '2.202'
$ws.Range('E24').Value = '  +0.12%  '

$ws.Range('D25').Value = '2.122.39'
$ws.Range('E25').Value = '  +1.31%  '

$ws.Range('D26').Value = '154.48'
$ws.Range('E26').Value = '  -1.35%  '

$ws.Range('D27').Value = '19.74'
$ws.Range('E27').Value = '  +1.21%  '

$ws.Range('D28').Value = '6.047'
$ws.Range('E28').Value = '  +10.69%  '

$ws.Range('D29').Value = '2.095'
$ws.Range('E29').Value = '  +0.85%  '

$ws.Range('D30').Value = '117.68'
$ws.Range('E30').Value = '  +0.25%  '

$ws.Range('D31').Value = '1.069'
$ws.Range('E31').Value = '  +4.65%  '

$ws.Range('D32').Value = '0.09473'
$ws.Range('E32').Value = '  +1.72%  '

$ws.Range('E33').Value = '  +0.75%  '

$ws.Range('D34').Value = '3.554'
$ws.Range('E34').Value = '  +1.03%  '

$ws.Range('D35').Value = '5.358'
$ws.Range('E35').Value = '  +1.54%  '

$ws.Range('D36').Value = '0.06081'
$ws.Range('E36').Value = '  +0.68%  '

$ws.Range('D37').Value = '0.02247'
$ws.Range('E37').Value = '  +0.98%  '

$ws.Range('D38').Value = '1.170'
$ws.Range('E38').Value = '  -0.55%  '

$ws.Range('D39').Value = '8.082'
$ws.Range('E39').Value = '  -2.70%  '

$ws.Range('D40').Value = '0.5813'
$ws.Range('E40').Value = '  +0.65%  '

$ws.Range('D41').Value = '2.516'
$ws.Range('E41').Value = '  +12.71%  '

$ws.Range('D42').Value = '0.1830'
$ws.Range('E42').Value = '  +0.43%  '

$ws.Range('D43').Value = '10.08'
$ws.Range('E43').Value = '  +0.32%  '

$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '1.272'
$ws.Range('E44').Value = '  +1.07%  '

$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '0.07754'
$ws.Range('E45').Value = '  +3.59%  '

$ws.Range('D46').Value = '12.13'
$ws.Range('E46').Value = '  +1.52%  '

$ws.Range('D47').Value = '0.5483'
$ws.Range('E47').Value = '  +0.87%  '

$ws.Range('D48').Value = '1.905'
$ws.Range('E48').Value = '  +0.41%  '

$ws.Range('D49').Value = '113.67'
$ws.Range('E49').Value = '  +2.30%  '

$ws.Range('D50').Value = '43.58'
$ws.Range('E50').Value = '  -3.47%  '

$ws.Range('D51').Value = '0.2914'
$ws.Range('E51').Value = '  +3.92%  '
